# Auto-generated edit script applying numeric updates to multiple sheets
# per the commit diff (Leve market-price / profit recalculation).
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H40").Value = 0
$ws.Range("I40").Value = 0
$ws.Range("J40").Value = 0
$ws.Range("K40").Value = 0
$ws.Range("L40").Value = 0
$ws.Range("M40").ClearContents()
$ws.Range("N40").ClearContents()

$ws.Range("H62").Value = 6500
$ws.Range("J62").Value = 6500
$ws.Range("L62").Value = 6500
$ws.Range("N62").Value = -7748

$ws.Range("H65").Value = 6500
$ws.Range("J65").Value = 6500
$ws.Range("L65").Value = 32500
$ws.Range("N65").Value = -38740

$ws.Range("H100").Value = 1433404.4
$ws.Range("I100").Value = 1671471.9
$ws.Range("K100").Value = 1671471.9
$ws.Range("M100").Value = -1670930.9

$ws.Range("H121").Value = 1949
$ws.Range("J121").Value = 1949
$ws.Range("L121").Value = 5847
$ws.Range("N121").Value = -9341

$ws.Range("H137").Value = 1250
$ws.Range("I137").Value = 1250
$ws.Range("K137").Value = 3750
$ws.Range("M137").Value = -1200

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 3218.889
$ws.Range("I32").Value = 3387.4783
$ws.Range("K32").Value = 3387.4783
$ws.Range("M32").Value = -3100.4783

$ws.Range("H61").Value = 6349.3335
$ws.Range("I61").Value = 6349.3335
$ws.Range("K61").Value = 6349.3335
$ws.Range("M61").Value = -6137.3335

$ws.Range("H63").Value = 0
$ws.Range("I63").Value = 0
$ws.Range("K63").Value = 0
$ws.Range("M63").ClearContents()

$ws.Range("H66").Value = 0
$ws.Range("I66").Value = 0
$ws.Range("K66").Value = 0
$ws.Range("M66").ClearContents()

$ws.Range("H104").Value = 0
$ws.Range("J104").Value = 0
$ws.Range("L104").Value = 0
$ws.Range("N104").ClearContents()

$ws.Range("H136").Value = 6349.3335
$ws.Range("I136").Value = 6349.3335
$ws.Range("K136").Value = 19048.0005
$ws.Range("M136").Value = -16498.0005

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H105").Value = 1990
$ws.Range("I105").Value = 1990
$ws.Range("K105").Value = 1990
$ws.Range("M105").Value = -243

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 1092.5555
$ws.Range("I31").Value = 1110.5
$ws.Range("K31").Value = 1110.5
$ws.Range("M31").Value = -815.5

$ws.Range("H34").Value = 1092.5555
$ws.Range("I34").Value = 1110.5
$ws.Range("K34").Value = 1110.5
$ws.Range("M34").Value = -908.5

$ws.Range("H43").Value = 18996.334
$ws.Range("J43").Value = 18996.334
$ws.Range("L43").Value = 18996.334
$ws.Range("N43").Value = -19364.334

$ws.Range("H58").Value = 2345
$ws.Range("J58").Value = 2345
$ws.Range("L58").Value = 2345
$ws.Range("N58").Value = -2751

$ws.Range("H95").Value = 56666.668
$ws.Range("J95").Value = 56666.668
$ws.Range("L95").Value = 56666.668
$ws.Range("N95").Value = -62158.668

$ws.Range("H97").Value = 37750
$ws.Range("J97").Value = 37750
$ws.Range("L97").Value = 37750
$ws.Range("N97").Value = -39732

$ws.Range("H101").Value = 18996.334
$ws.Range("J101").Value = 18996.334
$ws.Range("L101").Value = 18996.334
$ws.Range("N101").Value = -25486.334

$ws.Range("H102").Value = 0
$ws.Range("J102").Value = 0
$ws.Range("L102").Value = 0
$ws.Range("N102").ClearContents()

$ws.Range("H105").Value = 4504.75
$ws.Range("I105").Value = 6254.5
$ws.Range("K105").Value = 6254.5
$ws.Range("M105").Value = -4507.5

$ws.Range("H122").Value = 2026.2
$ws.Range("I122").Value = 1724.6666
$ws.Range("K122").Value = 5173.9998
$ws.Range("M122").Value = -2723.9998

$ws.Range("H132").Value = 10685.25
$ws.Range("I132").Value = 13497.9
$ws.Range("K132").Value = 40493.7
$ws.Range("M132").Value = -37963.7

$ws.Range("H136").Value = 2345
$ws.Range("J136").Value = 2345
$ws.Range("L136").Value = 7035
$ws.Range("N136").Value = -12135

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H3").Value = 0
$ws.Range("I3").Value = 0
$ws.Range("K3").Value = 0
$ws.Range("M3").ClearContents()

$ws.Range("H131").Value = 1557
$ws.Range("I131").Value = 1133
$ws.Range("J131").Value = 1875
$ws.Range("K131").Value = 3399
$ws.Range("L131").Value = 5625
$ws.Range("M131").Value = 1641
$ws.Range("N131").Value = -15705

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H141").Value = 0
$ws.Range("J141").Value = 0
$ws.Range("L141").Value = 0
$ws.Range("N141").ClearContents()

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H32").Value = 13
$ws.Range("I32").Value = 13
$ws.Range("K32").Value = 13
$ws.Range("M32").Value = 304

$ws.Range("H48").Value = 1599
$ws.Range("I48").Value = 1599
$ws.Range("K48").Value = 1599
$ws.Range("M48").Value = -938

$ws.Range("H64").Value = 150
$ws.Range("J64").Value = 150
$ws.Range("L64").Value = 150
$ws.Range("N64").Value = -600

$ws.Range("H67").Value = 150
$ws.Range("J67").Value = 150
$ws.Range("L67").Value = 150
$ws.Range("N67").Value = -1710

$ws.Range("H70").Value = 30000
$ws.Range("J70").Value = 30000
$ws.Range("L70").Value = 30000
$ws.Range("N70").Value = -30540

$ws.Range("H73").Value = 30000
$ws.Range("J73").Value = 30000
$ws.Range("L73").Value = 30000
$ws.Range("N73").Value = -31872

$ws.Range("H82").Value = 1474.8334
$ws.Range("I82").Value = 1999.5
$ws.Range("J82").Value = 1212.5
$ws.Range("K82").Value = 1999.5
$ws.Range("L82").Value = 1212.5
$ws.Range("M82").Value = -1638.5
$ws.Range("N82").Value = -1934.5

$ws.Range("H85").Value = 1474.8334
$ws.Range("I85").Value = 1999.5
$ws.Range("J85").Value = 1212.5
$ws.Range("K85").Value = 1999.5
$ws.Range("L85").Value = 1212.5
$ws.Range("M85").Value = -751.5
$ws.Range("N85").Value = -3708.5

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H64").Value = 25000
$ws.Range("I64").Value = 0
$ws.Range("J64").Value = 25000
$ws.Range("K64").Value = 0
$ws.Range("L64").Value = 25000
$ws.Range("N64").Value = -25496
$ws.Range("M64").ClearContents()

$ws.Range("H67").Value = 25000
$ws.Range("I67").Value = 0
$ws.Range("J67").Value = 25000
$ws.Range("K67").Value = 0
$ws.Range("L67").Value = 25000
$ws.Range("N67").Value = -26716
$ws.Range("M67").ClearContents()

$ws.Range("H70").Value = 39998
$ws.Range("J70").Value = 39997.5
$ws.Range("L70").Value = 39997.5
$ws.Range("N70").Value = -40627.5

$ws.Range("H73").Value = 39998
$ws.Range("J73").Value = 39997.5
$ws.Range("L73").Value = 39997.5
$ws.Range("N73").Value = -42181.5

$ws.Range("H76").Value = 30000
$ws.Range("J76").Value = 30000
$ws.Range("L76").Value = 30000
$ws.Range("N76").Value = -30630

$ws.Range("H79").Value = 30000
$ws.Range("J79").Value = 30000
$ws.Range("L79").Value = 30000
$ws.Range("N79").Value = -32184

$ws.Range("H126").Value = 1089.4286
$ws.Range("J126").Value = 2700
$ws.Range("L126").Value = 8100
$ws.Range("N126").Value = -13040
